$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - JAL
$ws.Range("C4").Value = 1688592
$ws.Range("N4").Value = 8.1
$ws.Range("P4").Value = 9
$ws.Range("R4").Value = 0.68
$ws.Range("V4").Value = 26

# Row 5 - JAL3
$ws.Range("C5").Value = 206218
$ws.Range("N5").Value = 8.1
$ws.Range("P5").Value = 9
$ws.Range("R5").Value = 0.68
$ws.Range("V5").Value = 26

# Row 6 - JFL
$ws.Range("C6").Value = 1226090
$ws.Range("N6").Value = 7.22
$ws.Range("P6").Value = 10
$ws.Range("R6").Value = 0.718
$ws.Range("V6").Value = 26

# Row 7 - JKL
$ws.Range("C7").Value = 1844988
$ws.Range("N7").Value = 8.16
$ws.Range("P7").Value = 8
$ws.Range("R7").Value = 0.603
$ws.Range("V7").Value = 26

# Row 8 - MFL
$ws.Range("C8").Value = 2255336
$ws.Range("N8").Value = 7.31
$ws.Range("R8").Value = 0.65
$ws.Range("V8").Value = 26

# Row 9 - FFL2
$ws.Range("C9").Value = 1365460
$ws.Range("N9").Value = 7.52
$ws.Range("P9").Value = 8
$ws.Range("V9").Value = 26

# Row 10 - JKL-U2
$ws.Range("C10").Value = 2667920
$ws.Range("N10").Value = 9.41
$ws.Range("P10").Value = 8
$ws.Range("R10").Value = 0.653
$ws.Range("V10").Value = 26

# Row 11 - GMT TOTAL (summary row inputs)
$ws.Range("N11").Value = 8.07
$ws.Range("P11").Value = 8.77
$ws.Range("R11").Value = 0.6589
$ws.Range("V11").Value = 26

# Row 12 - LINGERIE
$ws.Range("C12").Value = 3049596
$ws.Range("N12").Value = 4.13
$ws.Range("P12").Value = 10
$ws.Range("R12").Value = 0.73
$ws.Range("V12").Value = 26

# Row 13 - GTAL
$ws.Range("C13").Value = 241401
$ws.Range("N13").Value = 6.94
$ws.Range("P13").Value = 8
$ws.Range("R13").Value = 0.65
$ws.Range("V13").Value = 26

# Update selected cell to reflect the saved cursor position
$ws.Range("P14").Select()
